$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rule R30's "From" threshold (cell C10) from 18 to 1
$ws.Range("C10").Value = 1
